$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, pushing the existing Sprint 1 rows
# (old rows 6-11) down to rows 7-12, and leaving row 5 available (it was
# already blank-but-styled) to become a second "Sprint 2" entry.
$ws.Rows("6:6").Insert()

# Fill in the two new Sprint 2 task rows (row 5 was blank, row 6 is the
# newly inserted row) with the new task data.
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Oct 1 - Oct 6"
$ws.Range("D5").Value = "Write Sample Test Suite"

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Oct 1 - Oct 6"
$ws.Range("D6").Value = "Setup codeship and pipeline"

# Update the last selection, matching the saved file's view state.
$ws.Range("C19").Select()
